# Rename the second sheet ("Sheet2") to "adduser" and populate it with the
# new user-registration rows (username/password header + 3 sample rows),
# matching the upstream "issue in adding new excel sheet under
# src/main/resources" commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(2)
$ws.Name = "adduser"

# Make it the active sheet / tab (this also clears tabSelected on the
# previously-active "loginpage" sheet, matching the diff).
$ws.Activate()

$ws.Range("A1").Value = "usename"
$ws.Range("B1").Value = "password"

$ws.Range("A2").Value = "nee"
$ws.Range("B2").Value = 112323

$ws.Range("A3").Value = "neethu"
$ws.Range("B3").Value = 112323

$ws.Range("A4").Value = "nqww"
$ws.Range("B4").Value = 54546465

# Match the committed cursor position on the new sheet.
$ws.Range("D6").Select()
